{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The document's very last paragraph is an empty list placeholder that must\n// stay last, so anchor on the last paragraph that actually has text (i.e.\n// \"Apply css in child component\") and append the new bullets after it.\nconst items = paragraphs.items;\nlet anchor = null;\nfor (let i = items.length - 1; i >= 0; i--) {\n  if (items[i].text.trim().length > 0) {\n    anchor = items[i];\n    break;\n  }\n}\n\n// New bullet items: a \"Send data child to parent component\" topic, mirroring\n// the existing parent->child sub-list structure.\nconst newItems = [\n  { text: \"Send data from parent to child component\", level: 0 },\n  { text: \"Make child component\", level: 1 },\n  { text: \"Use child component in parent component\", level: 1 },\n  { text: \"Sending of a function from parent component to child\", level: 1 },\n  { text: \"Calling of that function from child \", level: 1 },\n  { text: \"Get data in parent component\", level: 1 },\n];\n\nfor (const item of newItems) {\n  const newPara = anchor.insertParagraph(item.text, \"After\");\n  const listItem = newPara.listItemOrNullObject;\n  listItem.level = item.level;\n  anchor = newPara;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# New bullet items to append under \"Reusable component\" topic, mirroring the\n# existing \"Send data from parent component\" sub-list but for child -> parent.\n$items = @(\n    @{ text = \"Send data from parent to child component\"; lvl = 1 },\n    @{ text = \"Make child component\"; lvl = 2 },\n    @{ text = \"Use child component in parent component\"; lvl = 2 },\n    @{ text = \"Sending of a function from parent component to child\"; lvl = 2 },\n    @{ text = \"Calling of that function from child \"; lvl = 2 },\n    @{ text = \"Get data in parent component\"; lvl = 2 }\n)\n\n$paras = $d.Paragraphs\n\n# Anchor on the last paragraph that actually contains text (the trailing\n# paragraph in the document body is an empty list placeholder that must stay\n# last), so the new bullets land right before it, e.g. after\n# \"Apply css in child component\".\n$anchorIndex = -1\nfor ($i = $paras.Count; $i -ge 1; $i--) {\n    $t = $paras.Item($i).Range.Text -replace \"[\\r\\a\\v]+$\", \"\"\n    if ($t.Trim().Length -gt 0) {\n        $anchorIndex = $i\n        break\n    }\n}\n\n$r = $paras.Item($anchorIndex).Range\n$r.Collapse(0)\n\nforeach ($item in $items) {\n    $r.InsertParagraphAfter()\n    $newp = $paras.Item($paras.Count - 1)\n    $newp.Range.Text = $item.text\n    $newp.Range.ListFormat.ListLevelNumber = $item.lvl\n    $r = $newp.Range\n    $r.Collapse(0)\n}\n"}
